# "First preparation for arena page was provided: common names and links, tabs, pairs list"
#
# This script:
#  1. Adds a new worksheet "location-page.html" at the end of the workbook
#     (becomes the active / selected tab).
#  2. Fills it with the origin/en/ua/ru translation rows for the new
#     "pairs list" arena screen.
#  3. Nudges the selection on the pre-existing sheets so that the header
#     row (A1:D12) is included, matching the reviewed state.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Create the new sheet after the last existing sheet.
# ---------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "location-page.html"

# ---------------------------------------------------------------------
# 2. Populate the translation table.
# ---------------------------------------------------------------------
$rows = @(
  @("origin", "en", "ua", "ru"),
  @("PHX_PAIRS_LIST", "Pairs list", "Список пар", "Список пар"),
  @("PHX_SETTINGS", "Settings", "Налаштування", "Настройки"),
  @("PHX_ROUND", "Round", "Раунд", "Раунд"),
  @("PHX_TIME", "Time", "Час", "Время"),
  @("PHX_RED_SPORTS", "Red corner", "Червоний кут", "Красный угол"),
  @("PHX_BLUE_SPORTS", "Blue corner", "Синій кут", "Синий угол"),
  @("PHX_GROUP", "Group", "Група", "Группа"),
  @("PHX_WINNER", "Winner", "Переможець", "Победитель"),
  @("PHX_RED_SCORE", "Red score", "Червоний рахунок", "Красный счёт"),
  @("PHX_BLUE_SCORE", "Blue score", "Синій рахунок", "Синий счёт"),
  @("PHX_PAIR_WINNER", "Winner of pair", "Переможець пари", "Победитель пари")
)

for ($r = 0; $r -lt $rows.Count; $r++) {
    $row = $rows[$r]
    for ($c = 0; $c -lt $row.Count; $c++) {
        $ws.Cells.Item($r + 1, $c + 1).Value = $row[$c]
    }
}

# ---------------------------------------------------------------------
# 3. Update the remembered selection on the other sheets so the new
#    common header block (A1:D12) is part of it.
# ---------------------------------------------------------------------
$common = $wb.Worksheets.Item("common")
$common.Range("A1:D12").Select()

$department = $wb.Worksheets.Item("department-page.html")
$department.Range("A1:D12").Select()

$trainer = $wb.Worksheets.Item("trainer-page.html")
$trainer.Range("A1:D12").Select()

$sportsman = $wb.Worksheets.Item("sportsman-page.html")
$sportsman.Range("A1:D12").Select()

$competition = $wb.Worksheets.Item("competition-page.html")
$competition.Range("A1:D12").Select()

$group = $wb.Worksheets.Item("group-page.html")
$group.Range("A1:D12").Select()

$pairsList = $wb.Worksheets.Item("creating-pairs-list-page.html")
$pairsList.Range("A1:D12").Select()

# ---------------------------------------------------------------------
# Make the new sheet the active tab again (it is the last one touched
# above, so re-select it last).
# ---------------------------------------------------------------------
$ws.Activate()
$ws.Range("A1:D12").Select()
